# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1) Insert a new "Player Info" sheet at the front with the player's bio.
# 2) Rename the MATCH_CARD_LINK column (on both existing sheets) to
#    MATCH_CODE and replace the full scorecard URL with just the numeric
#    match code that was embedded in it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "Player Info" sheet, inserted before "ODI Batting" so the sheet
#    order becomes: Player Info, ODI Batting, ODI Bowling.
#    (Do this first -- worksheet references shift position once a sheet
#    is inserted ahead of them, so grab the other sheets fresh afterward.)
# ---------------------------------------------------------------------
$battingForInsert = $wb.Worksheets.Item("ODI Batting")
$info = $wb.Worksheets.Add($battingForInsert)
$info.Name = "Player Info"

$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header formatting used by the
# header rows on the other sheets.
$infoHeader = $info.Range("A1:D1")
$infoHeader.Font.Bold = $true
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160
$infoHeader.Borders.LineStyle = 1
$infoHeader.Borders.Weight = 2

$info.Range("A2").NumberFormat = "@"
$info.Range("A2").Value = "4729"
$info.Range("B2").Value = "D'Arcy John Matthew Short"
$info.Range("C2").Value = "Left Handed"
$info.Range("D2").Value = "Left Arm Orthodox"

# ---------------------------------------------------------------------
# 2) ODI Batting: column D header + values (MATCH_CARD_LINK -> MATCH_CODE)
# ---------------------------------------------------------------------
$battingCodes = @{
    2 = "4167"
    3 = "4168"
    4 = "4170"
    5 = "4222"
    6 = "4415"
    7 = "4419"
    8 = "4421"
    9 = "4423"
}

$batting.Range("D1").Value = "MATCH_CODE"

foreach ($row in $battingCodes.Keys) {
    $cell = $batting.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$row]
}

# ---------------------------------------------------------------------
# 3) ODI Bowling: column B header + values (MATCH_CARD_LINK -> MATCH_CODE)
# ---------------------------------------------------------------------
$bowlingCodes = @{
    2 = "4167"
    3 = "4168"
    4 = "4170"
    5 = "4415"
    6 = "4421"
}

$bowling.Range("B1").Value = "MATCH_CODE"

foreach ($row in $bowlingCodes.Keys) {
    $cell = $bowling.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$row]
}

Write-Output "done"
